$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '72.115.46'
$ws.Range("E2").Value = "'" + '  +3.87%  '
$ws.Range("D3").Value = "'" + '3.624.76'
$ws.Range("E3").Value = "'" + '  +6.89%  '
$ws.Range("E4").Value = "'" + '  +0.04%  '
$ws.Range("D5").Value = "'" + '597.88'
$ws.Range("E5").Value = "'" + '  +1.79%  '
$ws.Range("D6").Value = "'" + '183.00'
$ws.Range("E6").Value = "'" + '  +1.53%  '
$ws.Range("D7").Value = "'" + '3.615.50'
$ws.Range("E7").Value = "'" + '  +6.87%  '
$ws.Range("D8").Value = "'" + '0.609'
$ws.Range("E8").Value = "'" + '  +2.00%  '
$ws.Range("E9").Value = "'" + '  +0.05%  '
$ws.Range("D10").Value = "'" + '0.207'
$ws.Range("E10").Value = "'" + '  +4.81%  '
$ws.Range("E11").Value = "'" + '  +2.99%  '
$ws.Range("D12").Value = "'" + '50.65'
$ws.Range("E12").Value = "'" + '  +4.23%  '
$ws.Range("E13").Value = "'" + '  +2.26%  '
$ws.Range("D14").Value = "'" + '701.29'
$ws.Range("E14").Value = "'" + '  +3.22%  '
$ws.Range("D15").Value = "'" + '4.205.14'
$ws.Range("E15").Value = "'" + '  +6.98%  '
$ws.Range("D16").Value = "'" + '8.97'
$ws.Range("E16").Value = "'" + '  +3.58%  '
$ws.Range("D17").Value = "'" + '72.148.27'
$ws.Range("E17").Value = "'" + '  +3.88%  '
$ws.Range("D18").Value = "'" + '3.602.60'
$ws.Range("E18").Value = "'" + '  +6.75%  '
$ws.Range("E19").Value = "'" + '  +1.72%  '
$ws.Range("D20").Value = "'" + '18.58'
$ws.Range("E20").Value = "'" + '  +4.83%  '
$ws.Range("E21").Value = "'" + '  +4.13%  '
$ws.Range("D22").Value = "'" + '0.936'
$ws.Range("E22").Value = "'" + '  +3.31%  '
$ws.Range("D23").Value = "'" + '5.84'
$ws.Range("E23").Value = "'" + '  +7.37%  '
$ws.Range("D24").Value = "'" + '17.93'
$ws.Range("E24").Value = "'" + '  +4.57%  '
$ws.Range("D25").Value = "'" + '105.56'
$ws.Range("E25").Value = "'" + '  +2.22%  '
$ws.Range("D26").Value = "'" + '4.03'
$ws.Range("E26").Value = "'" + '  +2.70%  '
$ws.Range("E27").Value = "'" + '  +4.48%  '
$ws.Range("E28").Value = "'" + '  +4.12%  '
$ws.Range("D29").Value = "'" + '35.84'
$ws.Range("E29").Value = "'" + '  +5.48%  '
$ws.Range("D30").Value = "'" + '9.14'
$ws.Range("E30").Value = "'" + '  +4.35%  '
$ws.Range("D31").Value = "'" + '7.43'
$ws.Range("E31").Value = "'" + '  +6.59%  '
$ws.Range("D32").Value = "'" + '4.21'
$ws.Range("E32").Value = "'" + '  +17.40%  '
$ws.Range("D33").Value = "'" + '592.88'
$ws.Range("E33").Value = "'" + '  +6.14%  '
$ws.Range("D34").Value = "'" + '11.38'
$ws.Range("E34").Value = "'" + '  +2.06%  '
$ws.Range("E35").Value = "'" + '  +1.54%  '
$ws.Range("D36").Value = "'" + '59.82'
$ws.Range("E36").Value = "'" + '  +2.06%  '
$ws.Range("E37").Value = "'" + '  +0.06%  '
$ws.Range("D38").Value = "'" + '0.146'
$ws.Range("E38").Value = "'" + '  +4.39%  '
$ws.Range("D39").Value = "'" + '3.651.87'
$ws.Range("E39").Value = "'" + '  -0.72%  '
$ws.Range("D40").Value = "'" + '36.06'
$ws.Range("E40").Value = "'" + '  +0.42%  '
$ws.Range("D41").Value = "'" + '0.0₃0778'
$ws.Range("E41").Value = "'" + '  +7.48%  '
$ws.Range("D42").Value = "'" + '3.47'
$ws.Range("E42").Value = "'" + '  +5.95%  '
$ws.Range("D43").Value = "'" + '2.79'
$ws.Range("E43").Value = "'" + '  +3.84%  '
$ws.Range("D44").Value = "'" + '0.0451'
$ws.Range("E44").Value = "'" + '  +6.68%  '
$ws.Range("E45").Value = "'" + '  +2.07%  '
$ws.Range("E46").Value = "'" + '  +2.18%  '
$ws.Range("E47").Value = "'" + '  +4.39%  '
$ws.Range("E48").Value = "'" + '  +4.59%  '
$ws.Range("E49").Value = "'" + '  +2.10%  '
$ws.Range("D50").Value = "'" + '0.998'
$ws.Range("E50").Value = "'" + '  -0.28%  '
$ws.Range("E51").Value = "'" + '  -0.01%  '
